$d = $word.ActiveDocument

function Get-ParaIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Table list: turn the blank paragraph right after "Shipment(...)" into the
#    new "Message(...)" table definition, and leave a fresh blank paragraph
#    behind it (before "FOREIGN KEY RELATIONSHIPS").
# ---------------------------------------------------------------------------
$idx = Get-ParaIndexByText("Shipment(ShipmentID, OrderID (FK), DateShipped, TrackingNum, Notes)")
$targetIdx = $idx + 1
$r = $d.Paragraphs($targetIdx).Range
$xml = $r.WordOpenXML

$oldPara = '<w:p w14:paraId="30AC4F23" w14:textId="77777777" w:rsidR="002B61C4" w:rsidRDefault="002B61C4" w:rsidP="00FF2BC8"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr></w:p>'
$newPara = '<w:p w14:paraId="30AC4F23" w14:textId="77777777" w:rsidR="002B61C4" w:rsidRDefault="002B61C4" w:rsidP="00FF2BC8"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Message(</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/><w:u w:val="single"/></w:rPr><w:t>MessageID</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve">EmployeeID(FK), </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve">CustomerID(FK), OrderID(FK), </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>MessageText, SentAt, SentBy)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr></w:p>'

$newXml = $xml.Replace($oldPara, $newPara)
$r.InsertXML($newXml)

# ---------------------------------------------------------------------------
# 2) Foreign-key list: turn the blank paragraph right after
#    "Shipment(OrderID) -> Order(OrderID)" into the three new Message(...)
#    relationship lines, leaving a fresh blank paragraph behind it (before
#    "COMPOSITE KEYS").
# ---------------------------------------------------------------------------
$idx2 = Get-ParaIndexByText("Shipment(OrderID) → Order(OrderID)")
$targetIdx2 = $idx2 + 1
$r2 = $d.Paragraphs($targetIdx2).Range
$xml2 = $r2.WordOpenXML

$oldPara2 = '<w:p w14:paraId="5AB9B798" w14:textId="77777777" w:rsidR="003E640B" w:rsidRDefault="003E640B" w:rsidP="003E640B"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr></w:p>'
$newPara2 = '<w:p w14:paraId="5AB9B798" w14:textId="77777777" w:rsidR="003E640B" w:rsidRDefault="003E640B" w:rsidP="003E640B"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Message(</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Employee</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>ID)</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>→</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Employee</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Employee</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>ID)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Message(</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Order</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>ID)</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>→</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Order(OrderID)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Message(CustomerID)</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>→</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> Customer(CustomerID)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr></w:p>'

$newXml2 = $xml2.Replace($oldPara2, $newPara2)
$r2.InsertXML($newXml2)

# ---------------------------------------------------------------------------
# 3) "COMPOSITE KEYS" heading: collapse the three runs ("COMPOSITE" / " KE" /
#    "YS") that make up the text into a single run. Round-tripping the
#    paragraph's own WordOpenXML merges adjacent runs that share identical
#    formatting, which is exactly what the diff shows.
# ---------------------------------------------------------------------------
$idx3 = Get-ParaIndexByText("COMPOSITE KEYS")
$r3 = $d.Paragraphs($idx3).Range
$r3.InsertXML($r3.WordOpenXML)
